$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting from the existing header cell (H1) onto the new headers
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for columns I (I0) and J (IF), rows 2-35
$data = @(
    @(2, 5),
    @(9, 9),
    @(7, 8),
    @(7, 8),
    @(9, 9),
    @(6, 8),
    @(4, 6),
    @(1, 3),
    @(1, 4),
    @(1, 5),
    @(1, 5),
    @(1, 3),
    @(1, 6),
    @(1, 7),
    @(1, 6),
    @(1, 5),
    @(1, 4),
    @(1, 5),
    @(1, 5),
    @(1, 5),
    @(3, 6),
    @(5, 7),
    @(6, 7),
    @(9, 9),
    @(9, 9),
    @(5, 7),
    @(6, 8),
    @(4, 5),
    @(7, 7),
    @(4, 5),
    @(5, 5),
    @(5, 6),
    @(6, 8),
    @(3, 4)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
